# Fill in the trial results for the Semi-Autonomous table (cols B-F, rows 3-9)
# and the Fully-Autonomous table (cols J-N, rows 3-14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Semi-Autonomous Data: rows 3..9 -> values for Trial 1..5 (columns B,C,D,E,F)
$semi = @(
    @(1,1,0,1,1),
    @(1,1,1,0,1),
    @(1,1,1,1,1),
    @(1,1,1,1,0),
    @(1,1,1,1,1),
    @(1,1,1,1,1),
    @(1,1,1,1,1)
)

for ($i = 0; $i -lt $semi.Length; $i++) {
    $row = 3 + $i
    for ($j = 0; $j -lt 5; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $semi[$i][$j]
    }
}

# Fully-Autonomous Data: rows 3..14 -> all successes (columns J,K,L,M,N)
for ($row = 3; $row -le 14; $row++) {
    for ($col = 10; $col -le 14; $col++) {
        $ws.Cells.Item($row, $col).Value = 1
    }
}

# Overall success-rate average under the Semi-Autonomous table
$ws.Range("G10").Formula = "=AVERAGE(G3:G9)"
$ws.Range("G10").NumberFormat = "0%"

$wb.Application.Calculate()

# New large title cell appended below the existing tables
$ws.Cells.Item(31, 9).Font.Size = 36
$ws.Cells.Item(31, 9).Font.Name = "Arial"
$ws.Cells.Item(31, 9).Font.Color = 2171169

$ws.Range("G24").Select()

$wb.Windows.Item(1).WindowState = -4143
